$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column I: header "Tempo veritiero" at I3, and "7gg" at I5
$ws.Range("I3").Value = "Tempo veritiero"
$ws.Range("I5").Value = "7gg"

# Set the column I width (closest reachable approximation to the target 15.28515625)
$ws.Columns.Item(9).ColumnWidth = 14.43

# Update selection to K5
$ws.Range("K5").Select()
